$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 150 - this shifts existing rows 150..231 down to 151..232,
# carrying their values/formatting with them (matching the cascading diff).
$ws.Rows(150).Insert()

# Populate the newly inserted row 150 with the new data point.
$ws.Range("A150").Value = 6
$ws.Range("B150").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C150").Value = "Metropolitana"
$ws.Range("D150").Value = 44806
$ws.Range("E150").Value = 13
$ws.Range("F150").Value = 100112029
$ws.Range("G150").Value = "Orégano"
$ws.Range("H150").Value = "Sin especificar"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 53
$ws.Range("K150").Value = 15000
$ws.Range("L150").Value = 16000
$ws.Range("M150").Value = 15472
$ws.Range("N150").Value = "$/docena de atados"
$ws.Range("O150").Value = "Región Metropolitana"
$ws.Range("P150").Value = 5157
$ws.Range("Q150").Value = 3
$ws.Range("R150").Value = "Hortaliza"
